$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be read as text so numeric-looking strings
# (e.g. "19.10", "1.000") keep their exact literal formatting
# instead of being auto-coerced to numbers by Excel.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "22.077.59"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "1.556.41"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "0.9991"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "290.44"
$ws.Range("E6").Value = "  +0.92%  "
$ws.Range("D7").Value = "0.3945"
$ws.Range("E7").Value = "  +3.92%  "
$ws.Range("D8").Value = "0.3241"
$ws.Range("E8").Value = "  -1.24%  "
$ws.Range("D9").Value = "43.31"
$ws.Range("E9").Value = "  -2.25%  "
$ws.Range("D10").Value = "0.07341"
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("D11").Value = "1.100"
$ws.Range("E11").Value = "  -4.53%  "
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "19.10"
$ws.Range("E13").Value = "  -7.09%  "
$ws.Range("D14").Value = "0.00001158"
$ws.Range("E14").Value = "  +6.20%  "
$ws.Range("D15").Value = "5.661"
$ws.Range("E15").Value = "  -3.43%  "
$ws.Range("D16").Value = "6.708"
$ws.Range("E16").Value = "  -1.68%  "
$ws.Range("D17").Value = "1.555.38"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").Value = "0.06611"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("D19").Value = "84.43"
$ws.Range("E19").Value = "  -2.14%  "
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").Value = "6.356"
$ws.Range("E21").Value = "  -0.91%  "
$ws.Range("D22").Value = "15.87"
$ws.Range("E22").Value = "  -2.59%  "
$ws.Range("D23").Value = "11.33"
$ws.Range("E23").Value = "  -3.45%  "
$ws.Range("D24").Value = "22.093.42"
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("D25").Value = "2.338"
$ws.Range("E25").Value = "  +1.59%  "
$ws.Range("D26").Value = "2.470"
$ws.Range("E26").Value = "  -4.32%  "
$ws.Range("D27").Value = "148.75"
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("D28").Value = "18.74"
$ws.Range("E28").Value = "  -3.64%  "
$ws.Range("D29").Value = "4.875"
$ws.Range("E29").Value = "  -1.48%  "
$ws.Range("D30").Value = "1.733.21"
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("D31").Value = "119.71"
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("D32").Value = "1.063"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("D33").Value = "5.742"
$ws.Range("E33").Value = "  -3.85%  "
$ws.Range("D34").Value = "0.08388"
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("D35").Value = "9.224"
$ws.Range("E35").Value = "  -3.98%  "
$ws.Range("D36").Value = "1.626"
$ws.Range("E36").Value = "  -14.36%  "
$ws.Range("D37").Value = "0.06253"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").Value = "0.02280"
$ws.Range("E38").Value = "  -5.09%  "
$ws.Range("D39").Value = "5.194"
$ws.Range("E39").Value = "  -2.99%  "
$ws.Range("E40").Value = "  -4.63%  "
$ws.Range("D41").Value = "1.216"
$ws.Range("E41").Value = "  -5.18%  "
$ws.Range("D42").Value = "10.80"
$ws.Range("E42").Value = "  -3.01%  "
$ws.Range("D43").Value = "0.9977"
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("D44").Value = "0.5873"
$ws.Range("E44").Value = "  -3.82%  "
$ws.Range("D45").Value = "13.18"
$ws.Range("E45").Value = "  -5.23%  "
$ws.Range("D46").Value = "3.761"
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("D47").Value = "0.5626"
$ws.Range("E47").Value = "  -5.49%  "
$ws.Range("D48").Value = "1.915"
$ws.Range("E48").Value = "  -4.54%  "
$ws.Range("D49").Value = "118.71"
$ws.Range("E49").Value = "  -4.51%  "
$ws.Range("D50").Value = "1.147"
$ws.Range("E50").Value = "  -2.62%  "
$ws.Range("D51").Value = "0.06871"
$ws.Range("E51").Value = "  -3.27%  "

# Restore the original (default/Normal) cell style on column D
# now that the text values are safely stored, so no visible
# formatting change is introduced.
$dRange.Style = "Normal"
